# Update countries & provincias Spain
# Refresh the COVID-19 "Pais" report table with newer figures and bump the
# "last updated" timestamp. A couple of countries leapfrogged each other in
# the ranking (Republica Dominicana overtook Panama; Libia overtook Hong
# Kong), so their name + stat rows trade places as well as get new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title row: refresh the "last updated" timestamp (16:35 -> 17:52)
$ws.Range("A1").Value = "Datos actualizados a 18 de Julio de 2020 a las 17:52"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 3788902
$ws.Range("C4").Value = 18890
$ws.Range("D4").Value = 1743987
$ws.Range("E4").Value = 1902536
$ws.Range("G4").Value = 315
$ws.Range("H4").Value = 142379

# India (row 6)
$ws.Range("B6").Value = 1055932
$ws.Range("C6").Value = 15475
$ws.Range("D6").Value = 664461
$ws.Range("E6").Value = 364963
$ws.Range("G6").Value = 223
$ws.Range("H6").Value = 26508

# row 13
$ws.Range("G13").Value = 40
$ws.Range("H13").Value = 45273

# row 17
$ws.Range("B17").Value = 244216
$ws.Range("C17").Value = 249
$ws.Range("D17").Value = 196806
$ws.Range("E17").Value = 12368
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 35042

# row 19
$ws.Range("B19").Value = 202442
$ws.Range("C19").Value = 97
$ws.Range("E19").Value = 5780

# Republica Dominicana now outranks Panama -> rows 42/43 swap names + stats
$ws.Range("A42").Value = "Republica Dominicana"
$ws.Range("B42").Value = 51519
$ws.Range("C42").Value = 1406
$ws.Range("D42").Value = 24607
$ws.Range("E42").Value = 25941
$ws.Range("G42").Value = 29
$ws.Range("H42").Value = 971

$ws.Range("A43").Value = "Panama"
$ws.Range("B43").Value = 51408
$ws.Range("D43").Value = 26520
$ws.Range("E43").Value = 23850
$ws.Range("H43").Value = 1038

# row 45
$ws.Range("D45").Value = 43833
$ws.Range("E45").Value = 3795

# row 101
$ws.Range("B101").Value = 3983
$ws.Range("C101").Value = 19
$ws.Range("E101").Value = 2415

# row 118
$ws.Range("B118").Value = 2072
$ws.Range("C118").Value = 48
$ws.Range("E118").Value = 1664
$ws.Range("H118").Value = 30

# Libia now outranks Hong Kong -> rows 126/127 swap names + stats
$ws.Range("A126").Value = "Libia"
$ws.Range("B126").Value = 1791
$ws.Range("C126").Value = 87
$ws.Range("D126").Value = 385
$ws.Range("E126").Value = 1358
$ws.Range("H126").Value = 48

$ws.Range("A127").Value = "Hong Kong"
$ws.Range("B127").Value = 1778
$ws.Range("C127").Value = 64
$ws.Range("D127").Value = 1274
$ws.Range("E127").Value = 492
$ws.Range("G127").Value = 1
$ws.Range("H127").Value = 12

# row 136
$ws.Range("B136").Value = 1348
$ws.Range("C136").Value = 12
$ws.Range("E136").Value = 203

# row 137
$ws.Range("B137").Value = 1214
$ws.Range("C137").Value = 5
$ws.Range("D137").Value = 1022
$ws.Range("E137").Value = 181
$ws.Range("G137").Value = 1
$ws.Range("H137").Value = 11
